$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.813.29'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.997.19'
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.89'
$ws.Range("E5").Value = '  +4.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.22'
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("E7").Value = '  +11.92%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.743'
$ws.Range("E9").Value = '  +1.13%  '
$ws.Range("E10").Value = '  -2.27%  '
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '50.27'
$ws.Range("E11").Value = '  +4.75%  '
$ws.Range("B12").Value = 'ShibaInu'
$ws.Range("C12").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000322'
$ws.Range("E12").Value = '  -3.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.67'
$ws.Range("E13").Value = '  -1.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.638.05'
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.001.69'
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.00'
$ws.Range("E16").Value = '  -0.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.42'
$ws.Range("E17").Value = '  -3.57%  '
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("E19").Value = '  -2.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.748.96'
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '426.90'
$ws.Range("E21").Value = '  -2.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '96.89'
$ws.Range("E22").Value = '  +0.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.48'
$ws.Range("E23").Value = '  -0.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.22'
$ws.Range("E24").Value = '  +5.66%  '
$ws.Range("E25").Value = '  -2.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.23'
$ws.Range("E26").Value = '  -5.45%  '
$ws.Range("E27").Value = '  -4.34%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.69'
$ws.Range("E28").Value = '  +19.34%  '
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.85'
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.71'
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.34'
$ws.Range("E31").Value = '  +4.98%  '
$ws.Range("E32").Value = '  +1.72%  '
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '673.62'
$ws.Range("E34").Value = '  -3.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '47.94'
$ws.Range("E35").Value = '  +18.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '64.92'
$ws.Range("E36").Value = '  -4.45%  '
$ws.Range("E37").Value = '  +2.07%  '
$ws.Range("E38").Value = '  -2.25%  '
$ws.Range("E39").Value = '  -8.75%  '
$ws.Range("E40").Value = '  -7.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.34'
$ws.Range("E42").Value = '  +5.62%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0486'
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("E45").Value = '  +3.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.83'
$ws.Range("E46").Value = '  +9.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.64'
$ws.Range("E47").Value = '  -3.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.34'
$ws.Range("E48").Value = '  -5.23%  '
$ws.Range("E49").Value = '  -3.90%  '
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.14'
$ws.Range("E51").Value = '  +0.31%  '
